# Auto-generated edit script applying cached-value updates from the commit diff.
# These H:N columns are plain cached numbers (no formulas) for currentAveragePrice /
# LevePrice / LeveProfit market-data columns, refreshed by the scheduled runner.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 112.92593   # H33: was 114.111115
$ws.Cells.Item(33, 9).Value = 119.333336   # I33: was 122
$ws.Cells.Item(33, 11).Value = 119.333336   # K33: was 122
$ws.Cells.Item(33, 13).Value = 109.666664   # M33: was 107
$ws.Cells.Item(98, 8).Value = 559547.2   # H98: was 621701.75
$ws.Cells.Item(98, 9).Value = 745164.5600000001   # I98: was 698664.7
$ws.Cells.Item(98, 10).Value = 2695.2   # J98: was 5998
$ws.Cells.Item(98, 11).Value = 745164.5600000001   # K98: was 698664.7
$ws.Cells.Item(98, 12).Value = 2695.2   # L98: was 5998
$ws.Cells.Item(98, 13).Value = -743666.5600000001   # M98: was -697166.7
$ws.Cells.Item(98, 14).Value = -5691.2   # N98: was -8994
$ws.Cells.Item(113, 8).Value = 2099.8   # H113: was 2349.75
$ws.Cells.Item(113, 9).Value = 2099.8   # I113: was 2349.75
$ws.Cells.Item(113, 11).Value = 2099.8   # K113: was 2349.75
$ws.Cells.Item(113, 13).Value = 1154.2   # M113: was 904.25
$ws.Cells.Item(118, 8).Value = 484   # H118: was 416.25
$ws.Cells.Item(118, 9).Value = 484   # I118: was 416.25
$ws.Cells.Item(118, 11).Value = 1452   # K118: was 1248.75
$ws.Cells.Item(118, 13).Value = 205   # M118: was 408.25
$ws.Cells.Item(122, 8).Value = 559547.2   # H122: was 621701.75
$ws.Cells.Item(122, 9).Value = 745164.5600000001   # I122: was 698664.7
$ws.Cells.Item(122, 10).Value = 2695.2   # J122: was 5998
$ws.Cells.Item(122, 11).Value = 2235493.68   # K122: was 2095994.1
$ws.Cells.Item(122, 12).Value = 8085.599999999999   # L122: was 17994
$ws.Cells.Item(122, 13).Value = -2233043.68   # M122: was -2093544.1
$ws.Cells.Item(122, 14).Value = -12985.6   # N122: was -22894

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 25225.256   # H32: was 25698.217
$ws.Cells.Item(32, 9).Value = 3036.5813   # I32: was 3026.2856
$ws.Cells.Item(32, 11).Value = 3036.5813   # K32: was 3026.2856
$ws.Cells.Item(32, 13).Value = -2749.5813   # M32: was -2739.2856
$ws.Cells.Item(122, 8).Value = 2079.35   # H122: was 1971.75
$ws.Cells.Item(122, 9).Value = 2038.2   # I122: was 2000.5625
$ws.Cells.Item(122, 10).Value = 2202.8   # J122: was 1914.125
$ws.Cells.Item(122, 11).Value = 6114.6   # K122: was 6001.6875
$ws.Cells.Item(122, 12).Value = 6608.400000000001   # L122: was 5742.375
$ws.Cells.Item(122, 13).Value = -3664.6   # M122: was -3551.6875
$ws.Cells.Item(122, 14).Value = -11508.4   # N122: was -10642.375
$ws.Cells.Item(132, 8).Value = 2516.9814   # H132: was 2895.8445
$ws.Cells.Item(132, 9).Value = 2066.587   # I132: was 2417.8108
$ws.Cells.Item(132, 11).Value = 6199.761   # K132: was 7253.432400000001
$ws.Cells.Item(132, 13).Value = -3669.761   # M132: was -4723.432400000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1376.2693   # H20: was 1475.125
$ws.Cells.Item(20, 9).Value = 1197.8235   # I20: was 1332.2
$ws.Cells.Item(20, 11).Value = 1197.8235   # K20: was 1332.2
$ws.Cells.Item(20, 13).Value = -950.8235   # M20: was -1085.2
$ws.Cells.Item(86, 8).Value = 6677.227   # H86: was 8379.8125
$ws.Cells.Item(86, 9).Value = 2549.3333   # I86: was 2507
$ws.Cells.Item(86, 10).Value = 11630.7   # J86: was 18167.834
$ws.Cells.Item(86, 11).Value = 2549.3333   # K86: was 2507
$ws.Cells.Item(86, 12).Value = 11630.7   # L86: was 18167.834
$ws.Cells.Item(86, 13).Value = -1426.3333   # M86: was -1384
$ws.Cells.Item(86, 14).Value = -13876.7   # N86: was -20413.834
$ws.Cells.Item(89, 8).Value = 6677.227   # H89: was 8379.8125
$ws.Cells.Item(89, 9).Value = 2549.3333   # I89: was 2507
$ws.Cells.Item(89, 10).Value = 11630.7   # J89: was 18167.834
$ws.Cells.Item(89, 11).Value = 12746.6665   # K89: was 12535
$ws.Cells.Item(89, 12).Value = 58153.5   # L89: was 90839.17
$ws.Cells.Item(89, 13).Value = -7130.666499999999   # M89: was -6919
$ws.Cells.Item(89, 14).Value = -69385.5   # N89: was -102071.17
$ws.Cells.Item(105, 8).Value = 297360.2   # H105: was 326063.1
$ws.Cells.Item(105, 9).Value = 3049.1667   # I105: was 3362.6316
$ws.Cells.Item(105, 10).Value = 1003706.6   # J105: was 837005.5
$ws.Cells.Item(105, 11).Value = 3049.1667   # K105: was 3362.6316
$ws.Cells.Item(105, 12).Value = 1003706.6   # L105: was 837005.5
$ws.Cells.Item(105, 13).Value = -1302.1667   # M105: was -1615.6316
$ws.Cells.Item(105, 14).Value = -1007200.6   # N105: was -840499.5
$ws.Cells.Item(134, 8).Value = 3086.5854   # H134: was 3317.5945
$ws.Cells.Item(134, 9).Value = 1981.3334   # I134: was 2160.739
$ws.Cells.Item(134, 11).Value = 5944.0002   # K134: was 6482.217000000001
$ws.Cells.Item(134, 13).Value = -3409.0002   # M134: was -3947.217000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1152   # H16: was 1071.125
$ws.Cells.Item(16, 9).Value = 1224.75   # I16: was 1090.6666
$ws.Cells.Item(16, 10).Value = 1006.5   # J16: was 1012.5
$ws.Cells.Item(16, 11).Value = 1224.75   # K16: was 1090.6666
$ws.Cells.Item(16, 12).Value = 1006.5   # L16: was 1012.5
$ws.Cells.Item(16, 13).Value = -937.75   # M16: was -803.6666
$ws.Cells.Item(16, 14).Value = -1580.5   # N16: was -1586.5
$ws.Cells.Item(58, 8).Value = 3627.7   # H58: was 2813.5652
$ws.Cells.Item(58, 9).Value = 2363.75   # I58: was 1128.9166
$ws.Cells.Item(58, 10).Value = 4470.3335   # J58: was 4651.364
$ws.Cells.Item(58, 11).Value = 2363.75   # K58: was 1128.9166
$ws.Cells.Item(58, 12).Value = 4470.3335   # L58: was 4651.364
$ws.Cells.Item(58, 13).Value = -2160.75   # M58: was -925.9166
$ws.Cells.Item(58, 14).Value = -4876.3335   # N58: was -5057.364
$ws.Cells.Item(92, 8).Value = 29266.334   # H92: was 29400
$ws.Cells.Item(92, 10).Value = 29266.334   # J92: was 29400
$ws.Cells.Item(92, 12).Value = 29266.334   # L92: was 29400
$ws.Cells.Item(92, 14).Value = -34258.334   # N92: was -34392
$ws.Cells.Item(99, 8).Value = 4816791.5   # H99: was 4816776
$ws.Cells.Item(99, 9).Value = 6259578.5   # I99: was 6259558.5
$ws.Cells.Item(99, 11).Value = 6259578.5   # K99: was 6259558.5
$ws.Cells.Item(99, 13).Value = -6258080.5   # M99: was -6258060.5
$ws.Cells.Item(100, 8).Value = 57390   # H100: was 68780
$ws.Cells.Item(100, 10).Value = 57390   # J100: was 68780
$ws.Cells.Item(100, 12).Value = 57390   # L100: was 68780
$ws.Cells.Item(100, 14).Value = -59554   # N100: was -70944
$ws.Cells.Item(113, 8).Value = 1152   # H113: was 1071.125
$ws.Cells.Item(113, 9).Value = 1224.75   # I113: was 1090.6666
$ws.Cells.Item(113, 10).Value = 1006.5   # J113: was 1012.5
$ws.Cells.Item(113, 11).Value = 1224.75   # K113: was 1090.6666
$ws.Cells.Item(113, 12).Value = 1006.5   # L113: was 1012.5
$ws.Cells.Item(113, 13).Value = 945.25   # M113: was 1079.3334
$ws.Cells.Item(113, 14).Value = -5346.5   # N113: was -5352.5
$ws.Cells.Item(126, 8).Value = 4816791.5   # H126: was 4816776
$ws.Cells.Item(126, 9).Value = 6259578.5   # I126: was 6259558.5
$ws.Cells.Item(126, 11).Value = 18778735.5   # K126: was 18778675.5
$ws.Cells.Item(126, 13).Value = -18776265.5   # M126: was -18776205.5
$ws.Cells.Item(136, 8).Value = 3627.7   # H136: was 2813.5652
$ws.Cells.Item(136, 9).Value = 2363.75   # I136: was 1128.9166
$ws.Cells.Item(136, 10).Value = 4470.3335   # J136: was 4651.364
$ws.Cells.Item(136, 11).Value = 7091.25   # K136: was 3386.7498
$ws.Cells.Item(136, 12).Value = 13411.0005   # L136: was 13954.092
$ws.Cells.Item(136, 13).Value = -4541.25   # M136: was -836.7498000000001
$ws.Cells.Item(136, 14).Value = -18511.0005   # N136: was -19054.092

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(20, 8).Value = 101   # H20: was 0
$ws.Cells.Item(20, 9).Value = 101   # I20: was 0
$ws.Cells.Item(20, 11).Value = 303   # K20: was 0
$ws.Cells.Item(20, 13).Value = -76   # M20: was None
$ws.Cells.Item(22, 8).Value = 654.5454999999999   # H22: was 900
$ws.Cells.Item(22, 10).Value = 1350   # J22: was 1800
$ws.Cells.Item(22, 12).Value = 4050   # L22: was 5400
$ws.Cells.Item(22, 14).Value = -4388   # N22: was -5738
$ws.Cells.Item(26, 8).Value = 253   # H26: was 204.75
$ws.Cells.Item(26, 9).Value = 221.875   # I26: was 191.5
$ws.Cells.Item(26, 10).Value = 502   # J26: was 271
$ws.Cells.Item(26, 11).Value = 665.625   # K26: was 574.5
$ws.Cells.Item(26, 12).Value = 1506   # L26: was 813
$ws.Cells.Item(26, 13).Value = -377.625   # M26: was -286.5
$ws.Cells.Item(26, 14).Value = -2082   # N26: was -1389
$ws.Cells.Item(27, 8).Value = 654.5454999999999   # H27: was 900
$ws.Cells.Item(27, 10).Value = 1350   # J27: was 1800
$ws.Cells.Item(27, 12).Value = 4050   # L27: was 5400
$ws.Cells.Item(27, 14).Value = -4254   # N27: was -5604
$ws.Cells.Item(32, 8).Value = 0   # H32: was 999
$ws.Cells.Item(32, 10).Value = 0   # J32: was 999
$ws.Cells.Item(32, 12).Value = 0   # L32: was 2997
$ws.Cells.Item(32, 14).Value = $null   # N32: clear (was -3563)
$ws.Cells.Item(34, 8).Value = 1428.2   # H34: was 1698.2
$ws.Cells.Item(34, 9).Value = 880.3333   # I34: was 1245.5
$ws.Cells.Item(34, 10).Value = 2250   # J34: was 2000
$ws.Cells.Item(34, 11).Value = 2640.9999   # K34: was 3736.5
$ws.Cells.Item(34, 12).Value = 6750   # L34: was 6000
$ws.Cells.Item(34, 13).Value = -2556.9999   # M34: was -3652.5
$ws.Cells.Item(34, 14).Value = -6918   # N34: was -6168
$ws.Cells.Item(39, 8).Value = 9389.474   # H39: was 9170
$ws.Cells.Item(39, 10).Value = 9389.474   # J39: was 9170
$ws.Cells.Item(39, 12).Value = 28168.422   # L39: was 27510
$ws.Cells.Item(39, 14).Value = -28756.422   # N39: was -28098
$ws.Cells.Item(46, 8).Value = 713.6   # H46: was 692.5333000000001
$ws.Cells.Item(46, 9).Value = 200   # I46: was 396
$ws.Cells.Item(46, 10).Value = 792.61536   # J46: was 766.6667
$ws.Cells.Item(46, 11).Value = 600   # K46: was 1188
$ws.Cells.Item(46, 12).Value = 2377.84608   # L46: was 2300.0001
$ws.Cells.Item(46, 13).Value = -509   # M46: was -1097
$ws.Cells.Item(46, 14).Value = -2559.84608   # N46: was -2482.0001
$ws.Cells.Item(56, 8).Value = 5153.636   # H56: was 5990
$ws.Cells.Item(56, 9).Value = 5153.636   # I56: was 5990
$ws.Cells.Item(56, 11).Value = 5153.636   # K56: was 5990
$ws.Cells.Item(56, 13).Value = -4623.636   # M56: was -5460
$ws.Cells.Item(58, 8).Value = 7395.316   # H58: was 7301.3
$ws.Cells.Item(58, 9).Value = 877.75   # I58: was 1805.2
$ws.Cells.Item(58, 11).Value = 2633.25   # K58: was 5415.6
$ws.Cells.Item(58, 13).Value = -2505.25   # M58: was -5287.6
$ws.Cells.Item(64, 8).Value = 2674.9092   # H64: was 1957.4615
$ws.Cells.Item(64, 9).Value = 736.6667   # I64: was 604.7778
$ws.Cells.Item(64, 10).Value = 5000.8   # J64: was 5001
$ws.Cells.Item(64, 11).Value = 2210.0001   # K64: was 1814.3334
$ws.Cells.Item(64, 12).Value = 15002.4   # L64: was 15003
$ws.Cells.Item(64, 13).Value = -1940.0001   # M64: was -1544.3334
$ws.Cells.Item(64, 14).Value = -15542.4   # N64: was -15543
$ws.Cells.Item(67, 8).Value = 2674.9092   # H67: was 1957.4615
$ws.Cells.Item(67, 9).Value = 736.6667   # I67: was 604.7778
$ws.Cells.Item(67, 10).Value = 5000.8   # J67: was 5001
$ws.Cells.Item(67, 11).Value = 2210.0001   # K67: was 1814.3334
$ws.Cells.Item(67, 12).Value = 15002.4   # L67: was 15003
$ws.Cells.Item(67, 13).Value = -1274.0001   # M67: was -878.3334
$ws.Cells.Item(67, 14).Value = -16874.4   # N67: was -16875
$ws.Cells.Item(70, 8).Value = 2560.7144   # H70: was 2600.1428
$ws.Cells.Item(70, 9).Value = 1106.25   # I70: was 1175.25
$ws.Cells.Item(70, 11).Value = 3318.75   # K70: was 3525.75
$ws.Cells.Item(70, 13).Value = -3003.75   # M70: was -3210.75
$ws.Cells.Item(73, 8).Value = 2560.7144   # H73: was 2600.1428
$ws.Cells.Item(73, 9).Value = 1106.25   # I73: was 1175.25
$ws.Cells.Item(73, 11).Value = 3318.75   # K73: was 3525.75
$ws.Cells.Item(73, 13).Value = -2226.75   # M73: was -2433.75
$ws.Cells.Item(75, 8).Value = 1754.9286   # H75: was 1849.55
$ws.Cells.Item(75, 9).Value = 909.75   # I75: was 1006.5
$ws.Cells.Item(75, 10).Value = 2093   # J75: was 2060.3125
$ws.Cells.Item(75, 11).Value = 2729.25   # K75: was 3019.5
$ws.Cells.Item(75, 12).Value = 6279   # L75: was 6180.9375
$ws.Cells.Item(75, 13).Value = -1731.25   # M75: was -2021.5
$ws.Cells.Item(75, 14).Value = -8275   # N75: was -8176.9375
$ws.Cells.Item(78, 8).Value = 1754.9286   # H78: was 1849.55
$ws.Cells.Item(78, 9).Value = 909.75   # I78: was 1006.5
$ws.Cells.Item(78, 10).Value = 2093   # J78: was 2060.3125
$ws.Cells.Item(78, 11).Value = 8187.75   # K78: was 9058.5
$ws.Cells.Item(78, 12).Value = 18837   # L78: was 18542.8125
$ws.Cells.Item(78, 13).Value = -3195.75   # M78: was -4066.5
$ws.Cells.Item(78, 14).Value = -28821   # N78: was -28526.8125
$ws.Cells.Item(82, 8).Value = 2085.9092   # H82: was 3822.2222
$ws.Cells.Item(82, 9).Value = 600   # I82: was 500
$ws.Cells.Item(82, 10).Value = 2234.5   # J82: was 4237.5
$ws.Cells.Item(82, 11).Value = 1800   # K82: was 1500
$ws.Cells.Item(82, 12).Value = 6703.5   # L82: was 12712.5
$ws.Cells.Item(82, 13).Value = -1394   # M82: was -1094
$ws.Cells.Item(82, 14).Value = -7515.5   # N82: was -13524.5
$ws.Cells.Item(85, 8).Value = 2085.9092   # H85: was 3822.2222
$ws.Cells.Item(85, 9).Value = 600   # I85: was 500
$ws.Cells.Item(85, 10).Value = 2234.5   # J85: was 4237.5
$ws.Cells.Item(85, 11).Value = 1800   # K85: was 1500
$ws.Cells.Item(85, 12).Value = 6703.5   # L85: was 12712.5
$ws.Cells.Item(85, 13).Value = -396   # M85: was -96
$ws.Cells.Item(85, 14).Value = -9511.5   # N85: was -15520.5
$ws.Cells.Item(94, 8).Value = 2801.6   # H94: was 5000
$ws.Cells.Item(94, 9).Value = 1012   # I94: was 0
$ws.Cells.Item(94, 10).Value = 3076.923   # J94: was 5000
$ws.Cells.Item(94, 11).Value = 3036   # K94: was 0
$ws.Cells.Item(94, 12).Value = 9230.769   # L94: was 15000
$ws.Cells.Item(94, 13).Value = -2360   # M94: was None
$ws.Cells.Item(94, 14).Value = -10582.769   # N94: was -16352
$ws.Cells.Item(104, 8).Value = 1900   # H104: was 0
$ws.Cells.Item(104, 9).Value = 1800   # I104: was 0
$ws.Cells.Item(104, 10).Value = 2000   # J104: was 0
$ws.Cells.Item(104, 11).Value = 5400   # K104: was 0
$ws.Cells.Item(104, 12).Value = 6000   # L104: was 0
$ws.Cells.Item(104, 13).Value = -2779   # M104: was None
$ws.Cells.Item(104, 14).Value = -11242   # N104: was None

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2352.0688   # H102: was 2437.7585
$ws.Cells.Item(102, 9).Value = 2399.7896   # I102: was 2530.5789
$ws.Cells.Item(102, 11).Value = 2399.7896   # K102: was 2530.5789
$ws.Cells.Item(102, 13).Value = -777.7896000000001   # M102: was -908.5789
$ws.Cells.Item(126, 8).Value = 2567.4666   # H126: was 2534.8333
$ws.Cells.Item(126, 9).Value = 2043.6666   # I126: was 2330.25
$ws.Cells.Item(126, 10).Value = 2916.6667   # J126: was 2616.6667
$ws.Cells.Item(126, 11).Value = 6130.9998   # K126: was 6990.75
$ws.Cells.Item(126, 12).Value = 8750.000100000001   # L126: was 7850.000100000001
$ws.Cells.Item(126, 13).Value = -3660.9998   # M126: was -4520.75
$ws.Cells.Item(126, 14).Value = -13690.0001   # N126: was -12790.0001
$ws.Cells.Item(137, 8).Value = 0   # H137: was 49780
$ws.Cells.Item(137, 10).Value = 0   # J137: was 49780
$ws.Cells.Item(137, 12).Value = 0   # L137: was 49780
$ws.Cells.Item(137, 14).Value = $null   # N137: clear (was -59980)

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 3500   # H40: was 3392.1428
$ws.Cells.Item(40, 9).Value = 0   # I40: was 2795
$ws.Cells.Item(40, 10).Value = 3500   # J40: was 3491.6667
$ws.Cells.Item(40, 11).Value = 0   # K40: was 2795
$ws.Cells.Item(40, 12).Value = 3500   # L40: was 3491.6667
$ws.Cells.Item(40, 13).Value = $null   # M40: clear (was -2659)
$ws.Cells.Item(40, 14).Value = -3772   # N40: was -3763.6667

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(110, 8).Value = 30395.2   # H110: was 32644
$ws.Cells.Item(110, 10).Value = 30395.2   # J110: was 32644
$ws.Cells.Item(110, 12).Value = 30395.2   # L110: was 32644
$ws.Cells.Item(110, 14).Value = -38575.2   # N110: was -40824
$ws.Cells.Item(126, 8).Value = 53349.95   # H126: was 63289.062
$ws.Cells.Item(126, 9).Value = 63108.688   # I126: was 72179.64
$ws.Cells.Item(126, 10).Value = 1303.3334   # J126: was 1055
$ws.Cells.Item(126, 11).Value = 189326.064   # K126: was 216538.92
$ws.Cells.Item(126, 12).Value = 3910.0002   # L126: was 3165
$ws.Cells.Item(126, 13).Value = -186856.064   # M126: was -214068.92
$ws.Cells.Item(126, 14).Value = -8850.0002   # N126: was -8105
$ws.Cells.Item(130, 8).Value = 48000   # H130: was 0
$ws.Cells.Item(130, 12).Value = 48000   # L130: was 0
$ws.Cells.Item(130, 14).Value = -58040   # N130: was None

Write-Output "Applied all cell updates"